# Applies the "Some Corrections on Form submission" edits to the
# Appointment Letter template: updated ref number, dates, candidate
# name/email/ID, designation, band/grade and the full CTC breakdown
# table (annual figure moved from 13,00,000 to 15,00,000).

$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $d.Content.Find.Execute(
        $findText,   # FindText
        $true,       # MatchCase
        $false,      # MatchWholeWord
        $false,      # MatchWildcards
        $false,      # MatchSoundsLike
        $false,      # MatchAllWordForms
        $true,       # Forward
        1,           # Wrap (wdFindContinue)
        $false,      # Format
        $replaceText,# ReplaceWith
        2            # Replace (wdReplaceAll)
    ) | Out-Null
}

# --- Band / Grade cells first -------------------------------------------
# These live in the CTC table next to a "Net Salary C=(A+B)" /
# "Total Compensation E=(C+D)" formula label, so a plain text search for
# the single letter "C" would also corrupt those formulas. Target the
# exact table cells instead.
$ctcTable = $d.Tables.Item(1)
$ctcTable.Rows.Item(3).Cells.Item(2).Range.Text = "B"    # Band
$ctcTable.Rows.Item(4).Cells.Item(2).Range.Text = "B2"   # Grade

# --- Reference number / dates --------------------------------------------
Replace-All "Appointment/25-26/17" "Appointment/25-26/19"
Replace-All "03-10-2025" "07-10-2025"
Replace-All "04-10-2025" "15-10-2025"

# --- Candidate details -----------------------------------------------------
Replace-All "Deepak Singh" "Swati Sharma"
Replace-All "Deepak@gmail.com" "swati.sharma@rigvedit.com"
Replace-All "111117" "111119"
Replace-All "Subject Matter Expert (SME)" "Senior Account Manager"

# --- Compensation summary ---------------------------------------------------
Replace-All "13,00,000" "15,00,000"
Replace-All "Thirteen Lakh Rupees Only" "Fifteen Lakh Rupees Only"

# --- CTC breakdown table (Annually / Monthly columns) -----------------------
Replace-All "4,16,000" "4,80,000"
Replace-All "34,667" "40,000"
Replace-All "2,08,000" "2,40,000"
Replace-All "17,333" "20,000"
Replace-All "41,600" "48,000"
Replace-All "3,467" "4,000"
Replace-All "49,920" "57,600"
Replace-All "4,160" "4,800"
Replace-All "2,45,936" "3,20,392"
Replace-All "20,495" "26,699"
Replace-All "11,44,856" "13,29,392"
Replace-All "95,405" "1,10,783"
Replace-All "60,256" "69,968"
Replace-All "5,021" "5,831"
Replace-All "12,05,112" "13,99,360"
Replace-All "1,00,426" "1,16,613"
Replace-All "19,968" "23,040"
Replace-All "1,664" "1,920"
Replace-All "25,000" "20,000"
Replace-All "2,083" "1,667"
Replace-All "94,888" "1,00,640"
Replace-All "7,907" "8,387"
Replace-All "1,08,333" "1,25,000"
Replace-All "52,420" "60,100"
Replace-All "4,360" "5,000"

Write-Host "Appointment letter corrections applied."
